# coursearrange: Add functions for course arrange by admin
#
# Adds a "课时" (class-hours) remark in column H of the "课程" (Course)
# sheet for a handful of rows, and makes the "课程" sheet the active /
# selected sheet (it previously was "教室").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("课程")

# Fill in the new remarks column (H) for the affected rows. The order in
# which the distinct text values are first written determines the order
# they are appended to the shared string table, so write them in the
# same sequence the original author did.
$ws.Range("H2").Value = "80课时"
$ws.Range("H3").Value = "80课时"
$ws.Range("H4").Value = "80课时"
$ws.Range("H13").Value = "4课时"
$ws.Range("H14").Value = "8课时"
$ws.Range("H12").Value = "60课时"
$ws.Range("H5").Value = "100课时"

# Make the "课程" sheet the active tab (previously "教室" was active),
# with H3 selected as the active cell.
$ws.Activate()
$ws.Range("H3").Select()
